$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the duplicate "comments" row (row 22). This is the extra comments
# section being removed per the commit message; the remaining "comments_2"
# field (now shifted up to row 24) is kept and turned into the single
# "comments" field, conditioned on the patient-tracing answer.
$ws.Rows.Item(22).Delete()

# H21 previously carried a one-off font style that existed only to support
# the row that just got removed; restore it to the sheet's normal style by
# copying the format from the neighboring cell that already uses it.
$ws.Cells.Item(21, 9).Copy()
$ws.Cells.Item(21, 8).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The row that used to hold "comments_2" (now row 24 after the delete above)
# becomes the single surviving "comments" field: rename it, give it the
# normal "Comments" label style (matching the style used elsewhere for that
# label), and make it conditional on the tracing question.
$ws.Cells.Item(24, 2).Value = "comments"

$ws.Cells.Item(22, 3).Copy()
$ws.Cells.Item(24, 3).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(24, 4).Value = "selected(${trace}, 'yes')"
